# Auto-generated cell updates derived from the OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    # Row 5
    "G5" = 1.5
    "H5" = 3.7
    "I5" = 6.7
    "M5" = 2.65
    "R5" = 2.2
    "S5" = 1.53
    "U5" = 5.9
    "W5" = 9.75
    "AA5" = 7.6
    "AB5" = 24
    "AC5" = 150
    "AE5" = 13.5
    # Row 6
    "G6" = 2.12
    "H6" = 3.1
    "I6" = 3.4
    "L6" = 1.5
    "M6" = 2.27
    "N6" = 2.4
    "O6" = 1.45
    "P6" = 1.52
    "Q6" = 2.2
    "R6" = 2.1
    "S6" = 1.57
    "U6" = 8.75
    "V6" = 9.5
    "W6" = 19.5
    "X6" = 22
    "Z6" = 6.6
    "AA6" = 6.2
    "AC6" = 150
    "AE6" = 7.6
    "AF6" = 16
    "AG6" = 13
    # Row 8
    "G8" = 1.7
    "H8" = 3.85
    "I8" = 4.1
    "L8" = 1.21
    "M8" = 3.55
    "N8" = 1.62
    "O8" = 2.02
    "R8" = 1.62
    "S8" = 2.02
    "T8" = 8.25
    "U8" = 9
    "W8" = 14
    "X8" = 12.5
    "Y8" = 22
    "Z8" = 13
    "AA8" = 7.7
    "AB8" = 14.5
    "AC8" = 55
    "AD8" = 400
    "AE8" = 13.5
    "AF8" = 24
    "AH8" = 65
    "AJ8" = 37
    # Row 9
    "G9" = 1.98
    "H9" = 3.65
    "I9" = 3.25
    "U9" = 11
    "W9" = 18.5
    "X9" = 14.5
    "Z9" = 13.5
    "AA9" = 7.3
    "AF9" = 19
    "AG9" = 11.25
    # Row 10
    "G10" = 1.57
    "H10" = 3.9
    "I10" = 5
    "L10" = 1.2
    "M10" = 3.6
    "N10" = 1.6
    "O10" = 2.05
    "R10" = 1.65
    "S10" = 1.98
    "T10" = 8.25
    "W10" = 12.5
    "X10" = 11.75
    "Y10" = 22
    "Z10" = 13
    "AA10" = 7.8
    "AB10" = 14.5
    "AC10" = 60
    "AD10" = 400
    "AE10" = 15.5
    "AF10" = 32
    "AG10" = 16
    "AH10" = 90
    "AJ10" = 45
    # Row 11
    "G11" = 4.6
    "H11" = 4.2
    "O11" = 2.42
    "T11" = 18.5
    "U11" = 32
    "V11" = 15
    "W11" = 80
    "X11" = 37
    "Z11" = 18
    "AB11" = 13.5
    "AE11" = 10.75
    "AF11" = 10
    "AH11" = 13.5
    "AI11" = 11.25
    # Row 13
    "G13" = 2.22
    "I13" = 3
    "P13" = 1.4
    "Q13" = 2.4
    "T13" = 5.7
    "U13" = 8.5
    "V13" = 7.8
    "W13" = 17.5
    "X13" = 16
    "Z13" = 7.8
    "AE13" = 7.1
    "AF13" = 12
    "AG13" = 9.25
    "AH13" = 29
    "AI13" = 22
    # Row 14
    "G14" = 3.6
    "I14" = 1.91
    "N14" = 2.1
    "P14" = 1.42
    "Q14" = 2.35
    "T14" = 7.5
    "U14" = 14.5
    "V14" = 10.75
    "W14" = 40
    "X14" = 29
    "Y14" = 40
    "AB14" = 14.5
    "AE14" = 5.2
    "AF14" = 6.9
    "AG14" = 7.4
    "AH14" = 13
    # Row 17
    "I17" = 2.8
    "N17" = 2.12
    # Row 18
    "T18" = 6.2
    "U18" = 8.5
    "X18" = 17
    "AD18" = 900
    "AE18" = 9.25
    "AF18" = 19
    "AI18" = 40
    # Row 24
    "G24" = 2.25
    "I24" = 2.75
    "N24" = 1.78
    "O24" = 1.83
    "T24" = 7.3
    "U24" = 9.75
    "V24" = 7.7
    "W24" = 18
    "X24" = 14.5
    "Y24" = 21
    "Z24" = 10.5
    "AD24" = 200
    "AE24" = 8.25
    "AF24" = 12
    "AG24" = 8.5
    "AH24" = 25
    "AI24" = 18
    "AJ24" = 23
    # Row 25
    "L25" = 1.53
    "M25" = 2.38
    # Row 32
    "G32" = 2.3
    "I32" = 3
    "J32" = 1.03
    "K32" = 10
    "P32" = 1.36
    "Q32" = 3
    "V32" = 9.5
    "W32" = 21
    "AG32" = 12
    # Row 34
    "G34" = 2.3
    "I34" = 3
    "R34" = 1.7
    "S34" = 2.05
    "T34" = 8.5
    "U34" = 12
    "V34" = 9.5
    "X34" = 19
    "AF34" = 15
    "AH34" = 29
    "AJ34" = 29
    # Row 35
    "G35" = 2.1
    "H35" = 3.7
    "I35" = 3.1
    "R35" = 1.5
    "S35" = 2.5
    "U35" = 13
    "W35" = 21
    "AD35" = 101
    "AE35" = 13
    "AF35" = 19
    "AG35" = 11
    "AH35" = 34
    "AI35" = 21
    "AJ35" = 23
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
